$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# R30 rule's "min" (From) value in column C, row 10, changes from 18 to 1.
$ws.Range("C10").Value = 1
